$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# ---------------------------------------------------------------------------
# 1) "Roll NO" / ":KWOWFL1165" currently live in two separate runs, split by
#    <w:proofErr> gramStart/gramEnd markers. Collapse them into one run and
#    drop the now-stale proofErr wrappers.
# ---------------------------------------------------------------------------
$rollFind = $d.Content
$rollFound = $rollFind.Find.Execute("Roll NO:KWOWFL1165")
if ($rollFound) {
    $rollPara = $rollFind.Paragraphs(1)
    $rollRange = $rollPara.Range
    $rollXml = '<w:p ' + $wNs + ' w:rsidR="00F77C8B" w:rsidRPr="00B93B4A" w:rsidRDefault="00F77C8B" w:rsidP="00F77C8B"><w:pPr><w:rPr><w:b/><w:color w:val="385623" w:themeColor="accent6" w:themeShade="80"/></w:rPr></w:pPr><w:r w:rsidRPr="00B93B4A"><w:rPr><w:b/><w:color w:val="385623" w:themeColor="accent6" w:themeShade="80"/></w:rPr><w:t>Roll NO:KWOWFL1165</w:t></w:r></w:p>'
    $rollRange.InsertXML($rollXml)
    Write-Host "Roll NO run merged"
} else {
    Write-Host "WARNING: 'Roll NO:KWOWFL1165' text not found"
}

# ---------------------------------------------------------------------------
# 2) The trailing paragraph ("  " + <w:bookmarkStart/bookmarkEnd name=_GoBack> +
#    "}") collapses into a single run "  }" and the _GoBack bookmark is
#    dropped from here (it gets re-created on the Assignment paragraph below).
#    Locate it via the bookmark itself so we don't depend on paragraph
#    indices.
# ---------------------------------------------------------------------------
$goBack = $d.Bookmarks("_GoBack")
$lastPara = $goBack.Range.Paragraphs(1)
$lastRange = $lastPara.Range
$lastXml = '<w:p ' + $wNs + ' w:rsidR="0021240D" w:rsidRPr="0021240D" w:rsidRDefault="0021240D"><w:pPr><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r w:rsidRPr="0021240D"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">  }</w:t></w:r></w:p>'
$lastRange.InsertXML($lastXml)
Write-Host "trailing bookmark paragraph merged"

# ---------------------------------------------------------------------------
# 3) "Assignment 3" -> "Assignment 4", and the _GoBack bookmark (now removed
#    from step 2) is re-added right after this run.
# ---------------------------------------------------------------------------
$asgFind = $d.Content
$asgFound = $asgFind.Find.Execute("Assignment 3")
if ($asgFound) {
    $asgPara = $asgFind.Paragraphs(1)
    $asgRange = $asgPara.Range
    $asgXml = '<w:p ' + $wNs + ' w:rsidR="00F77C8B" w:rsidRDefault="00F77C8B" w:rsidP="00F77C8B"><w:pPr><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="44"/><w:szCs w:val="44"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="44"/><w:szCs w:val="44"/></w:rPr><w:t>Assignment 4</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>'
    $asgRange.InsertXML($asgXml)
    Write-Host "Assignment 3 -> Assignment 4, bookmark relocated"
} else {
    Write-Host "WARNING: 'Assignment 3' text not found"
}

Write-Host "edits applied"
